$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing row, text + values updated, column B now used)
$ws.Range("A2").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_STB_BET"
$ws.Range("B2").Value = "BET"
$ws.Range("C2").Value = 45000
$ws.Range("D2").Value = 45003
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 3 (new row)
$ws.Range("A3").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_STB_BEW"
$ws.Range("B3").Value = "BEW"
$ws.Range("C3").Value = 44999
$ws.Range("D3").Value = 45001
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 4 (new row)
$ws.Range("A4").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_STB_Fertig"
$ws.Range("B4").Value = "Fertig"
$ws.Range("C4").Value = 44998
$ws.Range("D4").Value = 45006
$ws.Range("C4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 5 (new row)
$ws.Range("A5").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_STB_SCH"
$ws.Range("B5").Value = "SCH"
$ws.Range("C5").Value = 44997
$ws.Range("D5").Value = 44998
$ws.Range("C5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
